$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions bot) - updates Price (D) and Volume(1h) (E)
# columns for each coin row, plus a couple of rows that changed rank order.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.150.96"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.623.50"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.14"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  -1.83%  "

$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.31"
$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.626.47"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.122.49"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.55"
$ws.Range("E16").Value = "  -4.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0743"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.21"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.41"
$ws.Range("E22").Value = "  -6.89%  "

$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.16"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.29"
$ws.Range("E26").Value = "  -3.34%  "

$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0507"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("E30").Value = "  -1.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.36"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.341.22"
$ws.Range("E33").Value = "  +5.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  -0.66%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.554"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.73"
$ws.Range("E40").Value = "  +6.09%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.802"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.24"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.759.73"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.925"
$ws.Range("E45").Value = "  +37.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.30"
$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0997"
$ws.Range("E49").Value = "  +2.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("E51").Value = "  -0.13%  "
